# Auto-generated script to update Leve market-price columns (H:N) across sheets
# based on refreshed Universalis market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1311.1428
$ws.Range("I2").Value = 1256.8
$ws.Range("J2").Value = 1447
$ws.Range("K2").Value = 1256.8
$ws.Range("L2").Value = 1447
$ws.Range("M2").Value = -1143.8
$ws.Range("N2").Value = -1673
$ws.Range("H7").Value = 2501.6667
$ws.Range("I7").Value = 2002.5
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 2002.5
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -1890.5
$ws.Range("N7").Value = -3724
$ws.Range("H14").Value = 2501.6667
$ws.Range("I14").Value = 2002.5
$ws.Range("J14").Value = 3500
$ws.Range("K14").Value = 2002.5
$ws.Range("L14").Value = 3500
$ws.Range("M14").Value = -1811.5
$ws.Range("N14").Value = -3882
$ws.Range("H18").Value = 5776.5
$ws.Range("I18").Value = 5431.7
$ws.Range("J18").Value = 7500.5
$ws.Range("K18").Value = 5431.7
$ws.Range("L18").Value = 7500.5
$ws.Range("M18").Value = -5147.7
$ws.Range("N18").Value = -8068.5
$ws.Range("H20").Value = 4270
$ws.Range("I20").Value = 4270
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 4270
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -4040
$ws.Range("H32").Value = 3717.6
$ws.Range("I32").Value = 2149.5
$ws.Range("J32").Value = 9990
$ws.Range("K32").Value = 2149.5
$ws.Range("L32").Value = 9990
$ws.Range("M32").Value = -1823.5
$ws.Range("N32").Value = -10642
$ws.Range("H35").Value = 4270
$ws.Range("I35").Value = 4270
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 4270
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -3891
$ws.Range("H40").Value = 1308
$ws.Range("I40").Value = 1308
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1308
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1133
$ws.Range("H45").Value = 5444
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 5444
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 16332
$ws.Range("N45").Value = -16716
$ws.Range("H46").Value = 3332.6667
$ws.Range("I46").Value = 3998
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 11994
$ws.Range("L46").Value = 9000
$ws.Range("M46").Value = -11875
$ws.Range("N46").Value = -9238
$ws.Range("H48").Value = 11650
$ws.Range("I48").Value = 14950
$ws.Range("J48").Value = 10000
$ws.Range("K48").Value = 44850
$ws.Range("L48").Value = 30000
$ws.Range("M48").Value = -44558
$ws.Range("N48").Value = -30584
$ws.Range("H56").Value = 11650
$ws.Range("I56").Value = 14950
$ws.Range("J56").Value = 10000
$ws.Range("K56").Value = 44850
$ws.Range("L56").Value = 30000
$ws.Range("M56").Value = -44316
$ws.Range("N56").Value = -31068
$ws.Range("H59").Value = 2000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 2000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 6000
$ws.Range("N59").Value = -7114
$ws.Range("H60").Value = 3332.6667
$ws.Range("I60").Value = 3998
$ws.Range("J60").Value = 3000
$ws.Range("K60").Value = 11994
$ws.Range("L60").Value = 9000
$ws.Range("M60").Value = -11510
$ws.Range("N60").Value = -9968
$ws.Range("H132").Value = 6899.2
$ws.Range("I132").Value = 7332.4443
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 21997.3329
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -19467.3329
$ws.Range("M45").ClearContents()
$ws.Range("M59").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2310.3845
$ws.Range("I74").Value = 2294.75
$ws.Range("J74").Value = 2498
$ws.Range("K74").Value = 2294.75
$ws.Range("L74").Value = 2498
$ws.Range("M74").Value = -1420.75
$ws.Range("H77").Value = 2310.3845
$ws.Range("I77").Value = 2294.75
$ws.Range("J77").Value = 2498
$ws.Range("K77").Value = 11473.75
$ws.Range("L77").Value = 12490
$ws.Range("M77").Value = -7105.75
$ws.Range("H122").Value = 2097.8235
$ws.Range("I122").Value = 2436.077
$ws.Range("J122").Value = 998.5
$ws.Range("K122").Value = 7308.231000000001
$ws.Range("L122").Value = 2995.5
$ws.Range("M122").Value = -4858.231000000001
$ws.Range("N122").Value = -7895.5
$ws.Range("H125").Value = 90976.60000000001
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 90976.60000000001
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 90976.60000000001
$ws.Range("N125").Value = -100816.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 369.6
$ws.Range("I22").Value = 377.55554
$ws.Range("J22").Value = 298
$ws.Range("K22").Value = 377.55554
$ws.Range("L22").Value = 298
$ws.Range("M22").Value = -204.55554
$ws.Range("N22").Value = -644
$ws.Range("H94").Value = 516.55554
$ws.Range("I94").Value = 387.72726
$ws.Range("J94").Value = 719
$ws.Range("K94").Value = 387.72726
$ws.Range("L94").Value = 719
$ws.Range("M94").Value = 63.27274
$ws.Range("N94").Value = -1621
$ws.Range("H134").Value = 6280
$ws.Range("I134").Value = 6306.615
$ws.Range("J134").Value = 6107
$ws.Range("K134").Value = 18919.845
$ws.Range("L134").Value = 18321
$ws.Range("M134").Value = -16384.845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3636799.5
$ws.Range("I22").Value = 556.3333
$ws.Range("J22").Value = 5000390.5
$ws.Range("K22").Value = 556.3333
$ws.Range("L22").Value = 5000390.5
$ws.Range("M22").Value = -206.3333
$ws.Range("H31").Value = 2072.3333
$ws.Range("I31").Value = 2072.3333
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2072.3333
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1777.3333
$ws.Range("H34").Value = 2072.3333
$ws.Range("I34").Value = 2072.3333
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2072.3333
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1870.3333
$ws.Range("H132").Value = 3104.6365
$ws.Range("I132").Value = 3104.6365
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9313.9095
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6783.9095
$ws.Range("H134").Value = 3783.1428
$ws.Range("I134").Value = 3689.5386
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 11068.6158
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -8533.6158
$ws.Range("N134").Value = -20070
$ws.Range("N31").ClearContents()
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 3896.75
$ws.Range("I32").Value = 1799
$ws.Range("J32").Value = 4596
$ws.Range("K32").Value = 5397
$ws.Range("L32").Value = 13788
$ws.Range("M32").Value = -5114
$ws.Range("H92").Value = 220.33333
$ws.Range("I92").Value = 229.75
$ws.Range("J92").Value = 201.5
$ws.Range("K92").Value = 689.25
$ws.Range("L92").Value = 604.5
$ws.Range("M92").Value = 558.75
$ws.Range("N92").Value = -3100.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 706.5238000000001
$ws.Range("I22").Value = 457.46155
$ws.Range("J22").Value = 1111.25
$ws.Range("K22").Value = 457.46155
$ws.Range("L22").Value = 1111.25
$ws.Range("M22").Value = -162.46155
$ws.Range("H27").Value = 706.5238000000001
$ws.Range("I27").Value = 457.46155
$ws.Range("J27").Value = 1111.25
$ws.Range("K27").Value = 457.46155
$ws.Range("L27").Value = 1111.25
$ws.Range("M27").Value = -350.46155
$ws.Range("H46").Value = 2060.3635
$ws.Range("I46").Value = 1652.625
$ws.Range("J46").Value = 3147.6667
$ws.Range("K46").Value = 1652.625
$ws.Range("L46").Value = 3147.6667
$ws.Range("M46").Value = -1464.625
$ws.Range("N46").Value = -3523.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4671.8184
$ws.Range("I81").Value = 2673.75
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 5347.5
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = -4286.5
$ws.Range("H84").Value = 4671.8184
$ws.Range("I84").Value = 2673.75
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 26737.5
$ws.Range("L84").Value = 100000
$ws.Range("M84").Value = -21433.5
